$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Inscritos (column E) 22 -> 23
$ws.Range("E7").Value = 23

# Row 11: Inscritos (E) 10 -> 11, Pagos (F) 7 -> 8, Inscrições homologadas (H) 7 -> 8
$ws.Range("E11").Value = 11
$ws.Range("F11").Value = 8
$ws.Range("H11").Value = 8

# Row 16: Inscritos (column E) 268 -> 269
$ws.Range("E16").Value = 269
